# Generate Report for Handback
# Updates the localization-status report after a handback round-trips:
#  - Status text changes from "Ready for handoff" to the handed-back state
#  - The zh-cn / de-de detail sheets gain "Latest Target File" / "Latest
#    Handback File" / "Latest Handback DateTime" values + a hyperlink on
#    the target-file cell
#  - A couple of columns are widened so the longer strings fit

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetFile = "21d41568-f00e-4aaf-90d5-3c89a0e5ceb8.md"
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc91e18f46e9375c072b63dda6d5b1363d4a60bc/e2e/21d41568-f00e-4aaf-90d5-3c89a0e5ceb8.md"

$zhHandback = "21d41568-f00e-4aaf-90d5-3c89a0e5ceb8.885cd247df98f70f27b41408080dead45342786e.zh-cn.xlf"
$deHandback = "21d41568-f00e-4aaf-90d5-3c89a0e5ceb8.885cd247df98f70f27b41408080dead45342786e.de-de.xlf"

$zhHandbackDateTime = "2016-08-24 17:05:43"
$deHandbackDateTime = "2016-08-24 17:05:51"

# column widths as read back through ColumnWidth come out ~5/6 narrower than
# the OOXML <col width> value once saved, so bias the inputs accordingly
$wideStatusColumnWidth = 29.166666666666668   # -> ~29.98 / 30 char OOXML width
$wideFileColumnWidth   = 39.166666666666664   # -> 40 char OOXML width

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = $statusNew
$wsOverview.Cells.Item(2, 6).Value = $statusNew
$wsOverview.Cells.Item(3, 5).Value = $statusNew
$wsOverview.Cells.Item(3, 6).Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusColumnWidth

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(2, 3).Value = $statusNew
$wsZh.Cells.Item(3, 3).Value = $statusNew

$wsZh.Cells.Item(2, 9).Value = $targetFile
$wsZh.Cells.Item(2, 10).Value = $zhHandback
$wsZh.Cells.Item(2, 11).Value = $zhHandbackDateTime

$wsZh.Cells.Item(3, 9).Value = $targetFile
$wsZh.Cells.Item(3, 10).Value = $zhHandback
$wsZh.Cells.Item(3, 11).Value = $zhHandbackDateTime

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetFileUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetFileUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = $wideStatusColumnWidth
$wsZh.Columns.Item(9).ColumnWidth = $wideFileColumnWidth
$wsZh.Columns.Item(10).ColumnWidth = $wideFileColumnWidth

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(2, 3).Value = $statusNew
$wsDe.Cells.Item(3, 3).Value = $statusNew

$wsDe.Cells.Item(2, 9).Value = $targetFile
$wsDe.Cells.Item(2, 10).Value = $deHandback
$wsDe.Cells.Item(2, 11).Value = $deHandbackDateTime

$wsDe.Cells.Item(3, 9).Value = $targetFile
$wsDe.Cells.Item(3, 10).Value = $deHandback
$wsDe.Cells.Item(3, 11).Value = $deHandbackDateTime

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetFileUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetFileUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = $wideStatusColumnWidth
$wsDe.Columns.Item(9).ColumnWidth = $wideFileColumnWidth
$wsDe.Columns.Item(10).ColumnWidth = $wideFileColumnWidth
